# Apply "Add data for 2022-08-09" changes to the carjacking monthly workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet and update the "through" date references.
$ws.Name = "Through 2022-08-01"
$ws.Range("I1").Value = "2022 (through 08-01)"

# Update year-to-date totals for affected months.
$ws.Range("I3").Value = 140
$ws.Range("I8").Value = 167
$ws.Range("I9").Value = 2
$ws.Range("I14").Value = 974
